$d = $word.ActiveDocument
$d.Content.Find.Execute("letter_date_insert", $true, $false, $false, $false, $false, $true, 1, $false, "letter_date_insert", 2)
